# Auto-generated edit script applying odds updates from the commit diff
# "Atualizando o arquivo XLSX" - update numeric odds values across several rows
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H4").Value = 3
$ws.Range("I4").Value = 2.15
$ws.Range("K4").Value = 1.88
$ws.Range("L4").Value = 2.87
$ws.Range("W4").Value = 4.3
$ws.Range("AA4").Value = 2.12
$ws.Range("AH4").Value = 65
$ws.Range("AI4").Value = 6.2
$ws.Range("AJ4").Value = 6.1
$ws.Range("AK4").Value = 20
$ws.Range("AL4").Value = 150
$ws.Range("AM4").Value = 5.5
$ws.Range("AN4").Value = 8.75
$ws.Range("AP4").Value = 20
$ws.Range("AQ4").Value = 22
$ws.Range("J5").Value = 2.9
$ws.Range("L5").Value = 4.3
$ws.Range("AC5").Value = 5.6
$ws.Range("AD5").Value = 9
$ws.Range("AE5").Value = 9.5
$ws.Range("AG5").Value = 22
$ws.Range("AI5").Value = 6.2
$ws.Range("AM5").Value = 7.4
$ws.Range("AN5").Value = 16.5
$ws.Range("AO5").Value = 13.5
$ws.Range("AQ5").Value = 45
$ws.Range("AR5").Value = 65
$ws.Range("U6").Value = 3.95
$ws.Range("V6").Value = 1.24
$ws.Range("N7").Value = 8
$ws.Range("G8").Value = 3
$ws.Range("H8").Value = 2.55
$ws.Range("I8").Value = 2.82
$ws.Range("J8").Value = 3.75
$ws.Range("K8").Value = 1.78
$ws.Range("L8").Value = 3.55
$ws.Range("M8").Value = 1.16
$ws.Range("N8").Value = 4.55
$ws.Range("S8").Value = 2.85
$ws.Range("T8").Value = 1.37
$ws.Range("W8").Value = 5.2
$ws.Range("AC8").Value = 6.3
$ws.Range("AE8").Value = 11.5
$ws.Range("AG8").Value = 35
$ws.Range("AI8").Value = 4.55
$ws.Range("AJ8").Value = 5.2
$ws.Range("AK8").Value = 18
$ws.Range("AM8").Value = 6.1
$ws.Range("AN8").Value = 12.5
$ws.Range("AO8").Value = 11
$ws.Range("AP8").Value = 37
$ws.Range("AQ8").Value = 32
$ws.Range("M11").Value = 1.05
$ws.Range("N11").Value = 8
$ws.Range("O11").Value = 1.46
$ws.Range("P11").Value = 2.57
$ws.Range("Q11").Value = 1.83
$ws.Range("R11").Value = 2.03
$ws.Range("S11").Value = 2.38
$ws.Range("T11").Value = 1.57
$ws.Range("W11").Value = 4.5
$ws.Range("Y11").Value = 1.5
$ws.Range("Z11").Value = 2.37
$ws.Range("M12").Value = 1.04
$ws.Range("O12").Value = 1.27
$ws.Range("P12").Value = 3.5
$ws.Range("S12").Value = 2.03
$ws.Range("T12").Value = 1.83
$ws.Range("X12").Value = 1.27
$ws.Range("Y12").Value = 1.41
$ws.Range("Z12").Value = 2.62
$ws.Range("M13").Value = 1.07
$ws.Range("N13").Value = 7
$ws.Range("O13").Value = 1.49
$ws.Range("P13").Value = 2.45
$ws.Range("Q13").Value = 1.98
$ws.Range("R13").Value = 1.88
$ws.Range("S13").Value = 2.63
$ws.Range("T13").Value = 1.5
$ws.Range("W13").Value = 5
$ws.Range("X13").Value = 1.13
$ws.Range("Y13").Value = 1.58
$ws.Range("M14").Value = 1.05
$ws.Range("O14").Value = 1.37
$ws.Range("X14").Value = 1.19
$ws.Range("Y14").Value = 1.47
$ws.Range("AD14").Value = 11
$ws.Range("M15").Value = 1.02
$ws.Range("O15").Value = 1.15
$ws.Range("U15").Value = 2.05
$ws.Range("V15").Value = 1.8
$ws.Range("X15").Value = 1.47
$ws.Range("Y15").Value = 1.27
$ws.Range("G18").Value = 1.8
$ws.Range("H18").Value = 3.5
$ws.Range("I18").Value = 4.5
$ws.Range("AA18").Value = 1.8
$ws.Range("AB18").Value = 1.95
$ws.Range("AD18").Value = 8.5
$ws.Range("AM18").Value = 13
$ws.Range("AN18").Value = 23
$ws.Range("G19").Value = 2.5
$ws.Range("H19").Value = 3.1
$ws.Range("I19").Value = 2.8
$ws.Range("J19").Value = 3.25
$ws.Range("L19").Value = 3.6
$ws.Range("N19").Value = 9
$ws.Range("W19").Value = 4
$ws.Range("X19").Value = 1.22
$ws.Range("AC19").Value = 7.5
$ws.Range("AD19").Value = 12
$ws.Range("AE19").Value = 10
$ws.Range("AF19").Value = 23
$ws.Range("AM19").Value = 8.5
$ws.Range("AN19").Value = 13
$ws.Range("AO19").Value = 11
$ws.Range("AP19").Value = 29
$ws.Range("AQ19").Value = 23
$ws.Range("AR19").Value = 34
$ws.Range("G21").Value = 2.2
$ws.Range("I21").Value = 3.2
$ws.Range("J21").Value = 3
$ws.Range("T21").Value = 1.67
$ws.Range("AD21").Value = 10
$ws.Range("AE21").Value = 9.5
$ws.Range("AF21").Value = 21
$ws.Range("AH21").Value = 34
$ws.Range("AN21").Value = 15
$ws.Range("AP21").Value = 34
$ws.Range("T22").Value = 1.67
$ws.Range("T23").Value = 1.75
$ws.Range("N24").Value = 9
$ws.Range("S24").Value = 2.15
$ws.Range("S25").Value = 1.7
$ws.Range("S26").Value = 2.08
$ws.Range("T26").Value = 1.73
$ws.Range("K27").Value = 2.6
$ws.Range("L27").Value = 6.5
$ws.Range("M27").Value = 1.03
$ws.Range("N27").Value = 17
$ws.Range("O27").Value = 1.17
$ws.Range("P27").Value = 5
$ws.Range("S27").Value = 1.57
$ws.Range("T27").Value = 2.35
$ws.Range("U27").Value = 1.93
$ws.Range("V27").Value = 1.88
$ws.Range("W27").Value = 2.38
$ws.Range("X27").Value = 1.53
$ws.Range("AA27").Value = 1.8
$ws.Range("AB27").Value = 1.91
$ws.Range("AC27").Value = 8.5
$ws.Range("AD27").Value = 7.5
$ws.Range("AI27").Value = 17
$ws.Range("G29").Value = 1.47
$ws.Range("H29").Value = 4
$ws.Range("I29").Value = 6.1
$ws.Range("J29").Value = 2
$ws.Range("L29").Value = 5.8
$ws.Range("X29").Value = 1.36
$ws.Range("AA29").Value = 1.83
$ws.Range("AB29").Value = 1.78
$ws.Range("AC29").Value = 6.8
$ws.Range("AD29").Value = 6.9
$ws.Range("AF29").Value = 10.25
$ws.Range("AH29").Value = 26
$ws.Range("AI29").Value = 11.5
$ws.Range("AJ29").Value = 7.9
$ws.Range("AK29").Value = 17.5
$ws.Range("AM29").Value = 17
$ws.Range("AN29").Value = 40
$ws.Range("AO29").Value = 19
$ws.Range("I31").Value = 1.42
$ws.Range("S31").Value = 1.62
$ws.Range("AS31").Value = 600
$ws.Range("G32").Value = 3.75
$ws.Range("H32").Value = 3.6
$ws.Range("K32").Value = 2.2
$ws.Range("M32").Value = 1.02
$ws.Range("N32").Value = 12
$ws.Range("S32").Value = 1.77
$ws.Range("T32").Value = 1.97
$ws.Range("W32").Value = 3
$ws.Range("X32").Value = 1.36
$ws.Range("Y32").Value = 1.36
$ws.Range("Z32").Value = 3
$ws.Range("AA32").Value = 1.73
$ws.Range("AB32").Value = 2
$ws.Range("AC32").Value = 12
$ws.Range("AI32").Value = 12
$ws.Range("AJ32").Value = 7
$ws.Range("AM32").Value = 8
$ws.Range("AR32").Value = 26
$ws.Range("H33").Value = 3.25
$ws.Range("J33").Value = 3.75
$ws.Range("K33").Value = 2.1
$ws.Range("L33").Value = 3
$ws.Range("M33").Value = 1.06
$ws.Range("N33").Value = 9.5
$ws.Range("S33").Value = 2.03
$ws.Range("T33").Value = 1.75
$ws.Range("Y33").Value = 1.44
$ws.Range("Z33").Value = 2.63
$ws.Range("AA33").Value = 1.8
$ws.Range("AB33").Value = 1.91
$ws.Range("AC33").Value = 9.5
$ws.Range("AG33").Value = 26
$ws.Range("AI33").Value = 9.5
$ws.Range("AK33").Value = 15
$ws.Range("AL33").Value = 51
$ws.Range("AM33").Value = 7.5
$ws.Range("AQ33").Value = 19
$ws.Range("AR33").Value = 29
$ws.Range("AS33").Value = 251
